# "added max return portfolio calc"
# The "Dates & Risk Free Rate" sheet's End Date (B2) is rolled forward by one
# day (2023-05-19 -> 2023-05-20) to refresh the return-window used by the
# portfolio calculations, and the active selection is left on B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dates & Risk Free Rate")

$ws.Activate()
$ws.Range("B2").Value = 45066
$ws.Range("B2").Select()
